$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,2).Value = -17.36881335623137
$ws.Cells.Item(2,3).Value = 2.032184871412888
$ws.Cells.Item(2,4).Value = -17.36881335623137
$ws.Cells.Item(2,5).Value = -17.36881335623137
$ws.Cells.Item(2,6).Value = -17.36881335623137
$ws.Cells.Item(2,7).Value = -17.36881335623137
$ws.Cells.Item(2,8).Value = -17.36881335623137
$ws.Cells.Item(2,9).Value = -17.36881335623137
$ws.Cells.Item(2,10).Value = -17.36881335623137
$ws.Cells.Item(2,11).Value = -17.36881335623137
$ws.Cells.Item(3,2).Value = -17.36881335623137
$ws.Cells.Item(3,3).Value = -17.36881335623137
$ws.Cells.Item(3,4).Value = -17.36881335623137
$ws.Cells.Item(3,5).Value = -17.36881335623137
$ws.Cells.Item(3,6).Value = -17.36881335623137
$ws.Cells.Item(3,7).Value = -17.36881335623137
$ws.Cells.Item(3,8).Value = -17.36881335623137
$ws.Cells.Item(3,9).Value = 3.094484500542969
$ws.Cells.Item(3,10).Value = -17.36881335623137
$ws.Cells.Item(3,11).Value = -17.36881335623137
$ws.Cells.Item(4,2).Value = -17.36881335623137
$ws.Cells.Item(4,3).Value = 2.228840966238166
$ws.Cells.Item(4,4).Value = 2.183775070691092
$ws.Cells.Item(4,5).Value = -17.36881335623137
$ws.Cells.Item(4,6).Value = 3.385722697984496
$ws.Cells.Item(4,7).Value = -17.36881335623137
$ws.Cells.Item(4,8).Value = -17.36881335623137
$ws.Cells.Item(4,9).Value = -17.36881335623137
$ws.Cells.Item(4,10).Value = 2.210725406894172
$ws.Cells.Item(4,11).Value = -17.36881335623137
$ws.Cells.Item(5,2).Value = -17.36881335623137
$ws.Cells.Item(5,3).Value = 1.926365901715677
$ws.Cells.Item(5,4).Value = -17.36881335623137
$ws.Cells.Item(5,5).Value = -17.36881335623137
$ws.Cells.Item(5,6).Value = -17.36881335623137
$ws.Cells.Item(5,7).Value = 2.968138790357687
$ws.Cells.Item(5,8).Value = -17.36881335623137
$ws.Cells.Item(5,9).Value = -17.36881335623137
$ws.Cells.Item(5,10).Value = -17.36881335623137
$ws.Cells.Item(5,11).Value = -17.36881335623137
$ws.Cells.Item(6,2).Value = -17.36881335623137
$ws.Cells.Item(6,3).Value = -17.36881335623137
$ws.Cells.Item(6,4).Value = -17.36881335623137
$ws.Cells.Item(6,5).Value = -17.36881335623137
$ws.Cells.Item(6,6).Value = -17.36881335623137
$ws.Cells.Item(6,7).Value = -17.36881335623137
$ws.Cells.Item(6,8).Value = -17.36881335623137
$ws.Cells.Item(6,9).Value = -17.36881335623137
$ws.Cells.Item(6,10).Value = -17.36881335623137
$ws.Cells.Item(6,11).Value = -17.36881335623137
$ws.Cells.Item(7,2).Value = 2.699210434856341
$ws.Cells.Item(7,3).Value = -17.36881335623137
$ws.Cells.Item(7,4).Value = -17.36881335623137
$ws.Cells.Item(7,5).Value = -17.36881335623137
$ws.Cells.Item(7,6).Value = -17.36881335623137
$ws.Cells.Item(7,7).Value = -17.36881335623137
$ws.Cells.Item(7,8).Value = -17.36881335623137
$ws.Cells.Item(7,9).Value = -17.36881335623137
$ws.Cells.Item(7,10).Value = -17.36881335623137
$ws.Cells.Item(7,11).Value = -17.36881335623137
$ws.Cells.Item(8,2).Value = -17.36881335623137
$ws.Cells.Item(8,3).Value = -17.36881335623137
$ws.Cells.Item(8,4).Value = -17.36881335623137
$ws.Cells.Item(8,5).Value = 1.731759160550292
$ws.Cells.Item(8,6).Value = -17.36881335623137
$ws.Cells.Item(8,7).Value = -17.36881335623137
$ws.Cells.Item(8,8).Value = -17.36881335623137
$ws.Cells.Item(8,9).Value = -17.36881335623137
$ws.Cells.Item(8,10).Value = -17.36881335623137
$ws.Cells.Item(8,11).Value = -17.36881335623137
$ws.Cells.Item(9,2).Value = 3.755467645287628
$ws.Cells.Item(9,3).Value = -17.36881335623137
$ws.Cells.Item(9,4).Value = -17.36881335623137
$ws.Cells.Item(9,5).Value = -17.36881335623137
$ws.Cells.Item(9,6).Value = -17.36881335623137
$ws.Cells.Item(9,7).Value = -17.36881335623137
$ws.Cells.Item(9,8).Value = -17.36881335623137
$ws.Cells.Item(9,9).Value = -17.36881335623137
$ws.Cells.Item(9,10).Value = -17.36881335623137
$ws.Cells.Item(9,11).Value = -17.36881335623137
$ws.Cells.Item(10,2).Value = -17.36881335623137
$ws.Cells.Item(10,3).Value = -17.36881335623137
$ws.Cells.Item(10,4).Value = -17.36881335623137
$ws.Cells.Item(10,5).Value = -17.36881335623137
$ws.Cells.Item(10,6).Value = -17.36881335623137
$ws.Cells.Item(10,7).Value = -17.36881335623137
$ws.Cells.Item(10,8).Value = -17.36881335623137
$ws.Cells.Item(10,9).Value = 1.49409334350251
$ws.Cells.Item(10,10).Value = -17.36881335623137
$ws.Cells.Item(10,11).Value = 1.998492605182258
$ws.Cells.Item(11,2).Value = -17.36881335623137
$ws.Cells.Item(11,3).Value = -17.36881335623137
$ws.Cells.Item(11,4).Value = -17.36881335623137
$ws.Cells.Item(11,5).Value = 2.829920626950158
$ws.Cells.Item(11,6).Value = -17.36881335623137
$ws.Cells.Item(11,7).Value = 2.611439206520014
$ws.Cells.Item(11,8).Value = -17.36881335623137
$ws.Cells.Item(11,9).Value = -17.36881335623137
$ws.Cells.Item(11,10).Value = -17.36881335623137
$ws.Cells.Item(11,11).Value = 1.605402351848167
$ws.Cells.Item(12,2).Value = -17.36881335623137
$ws.Cells.Item(12,3).Value = -17.36881335623137
$ws.Cells.Item(12,4).Value = -17.36881335623137
$ws.Cells.Item(12,5).Value = -17.36881335623137
$ws.Cells.Item(12,6).Value = -17.36881335623137
$ws.Cells.Item(12,7).Value = -17.36881335623137
$ws.Cells.Item(12,8).Value = -17.36881335623137
$ws.Cells.Item(12,9).Value = -17.36881335623137
$ws.Cells.Item(12,10).Value = -17.36881335623137
$ws.Cells.Item(12,11).Value = -17.36881335623137
$ws.Cells.Item(13,2).Value = -17.36881335623137
$ws.Cells.Item(13,3).Value = -17.36881335623137
$ws.Cells.Item(13,4).Value = -17.36881335623137
$ws.Cells.Item(13,5).Value = 2.396642721098931
$ws.Cells.Item(13,6).Value = -17.36881335623137
$ws.Cells.Item(13,7).Value = -17.36881335623137
$ws.Cells.Item(13,8).Value = -17.36881335623137
$ws.Cells.Item(13,9).Value = -17.36881335623137
$ws.Cells.Item(13,10).Value = 2.273211886856974
$ws.Cells.Item(13,11).Value = 1.903490652097937
$ws.Cells.Item(14,2).Value = -17.36881335623137
$ws.Cells.Item(14,3).Value = -17.36881335623137
$ws.Cells.Item(14,4).Value = 1.254757396332281
$ws.Cells.Item(14,5).Value = -17.36881335623137
$ws.Cells.Item(14,6).Value = -17.36881335623137
$ws.Cells.Item(14,7).Value = -17.36881335623137
$ws.Cells.Item(14,8).Value = -17.36881335623137
$ws.Cells.Item(14,9).Value = -17.36881335623137
$ws.Cells.Item(14,10).Value = -17.36881335623137
$ws.Cells.Item(14,11).Value = 2.132870300959873
$ws.Cells.Item(15,2).Value = -17.36881335623137
$ws.Cells.Item(15,3).Value = -17.36881335623137
$ws.Cells.Item(15,4).Value = 1.30826276972302
$ws.Cells.Item(15,5).Value = -17.36881335623137
$ws.Cells.Item(15,6).Value = -17.36881335623137
$ws.Cells.Item(15,7).Value = -17.36881335623137
$ws.Cells.Item(15,8).Value = -17.36881335623137
$ws.Cells.Item(15,9).Value = -17.36881335623137
$ws.Cells.Item(15,10).Value = -17.36881335623137
$ws.Cells.Item(15,11).Value = -17.36881335623137
$ws.Cells.Item(16,2).Value = -17.36881335623137
$ws.Cells.Item(16,3).Value = -17.36881335623137
$ws.Cells.Item(16,4).Value = -17.36881335623137
$ws.Cells.Item(16,5).Value = -17.36881335623137
$ws.Cells.Item(16,6).Value = -17.36881335623137
$ws.Cells.Item(16,7).Value = -17.36881335623137
$ws.Cells.Item(16,8).Value = -17.36881335623137
$ws.Cells.Item(16,9).Value = -17.36881335623137
$ws.Cells.Item(16,10).Value = 2.434884218955572
$ws.Cells.Item(16,11).Value = -17.36881335623137
$ws.Cells.Item(17,2).Value = -17.36881335623137
$ws.Cells.Item(17,3).Value = 1.869556393881108
$ws.Cells.Item(17,4).Value = 2.223156329999242
$ws.Cells.Item(17,5).Value = -17.36881335623137
$ws.Cells.Item(17,6).Value = -17.36881335623137
$ws.Cells.Item(17,7).Value = -17.36881335623137
$ws.Cells.Item(17,8).Value = -17.36881335623137
$ws.Cells.Item(17,9).Value = 1.352796545213787
$ws.Cells.Item(17,10).Value = 1.497840754598765
$ws.Cells.Item(17,11).Value = -17.36881335623137
$ws.Cells.Item(18,2).Value = -17.36881335623137
$ws.Cells.Item(18,3).Value = -17.36881335623137
$ws.Cells.Item(18,4).Value = -17.36881335623137
$ws.Cells.Item(18,5).Value = -17.36881335623137
$ws.Cells.Item(18,6).Value = -17.36881335623137
$ws.Cells.Item(18,7).Value = -17.36881335623137
$ws.Cells.Item(18,8).Value = -17.36881335623137
$ws.Cells.Item(18,9).Value = 0.8975568638459454
$ws.Cells.Item(18,10).Value = 1.205044593469629
$ws.Cells.Item(18,11).Value = -17.36881335623137
$ws.Cells.Item(19,2).Value = -17.36881335623137
$ws.Cells.Item(19,3).Value = -17.36881335623137
$ws.Cells.Item(19,4).Value = 1.665483210184716
$ws.Cells.Item(19,5).Value = -17.36881335623137
$ws.Cells.Item(19,6).Value = -17.36881335623137
$ws.Cells.Item(19,7).Value = -17.36881335623137
$ws.Cells.Item(19,8).Value = -17.36881335623137
$ws.Cells.Item(19,9).Value = 1.500675061745104
$ws.Cells.Item(19,10).Value = -17.36881335623137
$ws.Cells.Item(19,11).Value = -17.36881335623137
$ws.Cells.Item(20,2).Value = -17.36881335623137
$ws.Cells.Item(20,3).Value = 0.8413870944914736
$ws.Cells.Item(20,4).Value = 1.46072667710926
$ws.Cells.Item(20,5).Value = -17.36881335623137
$ws.Cells.Item(20,6).Value = 3.255165448702976
$ws.Cells.Item(20,7).Value = -17.36881335623137
$ws.Cells.Item(20,8).Value = 4.321919997115192
$ws.Cells.Item(20,9).Value = 0.4801090551180914
$ws.Cells.Item(20,10).Value = -17.36881335623137
$ws.Cells.Item(20,11).Value = 2.273282418685074
$ws.Cells.Item(21,2).Value = -17.36881335623137
$ws.Cells.Item(21,3).Value = 0.9819747887242679
$ws.Cells.Item(21,4).Value = -17.36881335623137
$ws.Cells.Item(21,5).Value = 2.105163840166181
$ws.Cells.Item(21,6).Value = -17.36881335623137
$ws.Cells.Item(21,7).Value = 2.600152898735733
$ws.Cells.Item(21,8).Value = -17.36881335623137
$ws.Cells.Item(21,9).Value = -17.36881335623137
$ws.Cells.Item(21,10).Value = -17.36881335623137
$ws.Cells.Item(21,11).Value = -17.36881335623137
